$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 52 (pushes existing rows 52.. down to 54..)
$ws.Rows("52:53").Insert()

# Fill in new row 52 (Primera, week of 2022-12-20)
$ws.Cells.Item(52,1).Value = 11
$ws.Cells.Item(52,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(52,3).Value = "Bíobío"
$ws.Cells.Item(52,4).Value = 44915
$ws.Cells.Item(52,5).Value = 8
$ws.Cells.Item(52,6).Value = "Fruta"
$ws.Cells.Item(52,7).Value = 100101
$ws.Cells.Item(52,8).Value = "Berries"
$ws.Cells.Item(52,9).Value = 100101001
$ws.Cells.Item(52,10).Value = "Arándano (blue)"
$ws.Cells.Item(52,11).Value = "Sin especificar"
$ws.Cells.Item(52,12).Value = "Primera"
$ws.Cells.Item(52,13).Value = 200
$ws.Cells.Item(52,14).Value = 3200
$ws.Cells.Item(52,15).Value = 3500
$ws.Cells.Item(52,16).Value = 3350
$ws.Cells.Item(52,17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(52,18).Value = "Región de Ñuble"
$ws.Cells.Item(52,19).Value = 1675
$ws.Cells.Item(52,20).Value = 2

# Fill in new row 53 (Segunda, week of 2022-12-20)
$ws.Cells.Item(53,1).Value = 11
$ws.Cells.Item(53,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(53,3).Value = "Bíobío"
$ws.Cells.Item(53,4).Value = 44915
$ws.Cells.Item(53,5).Value = 8
$ws.Cells.Item(53,6).Value = "Fruta"
$ws.Cells.Item(53,7).Value = 100101
$ws.Cells.Item(53,8).Value = "Berries"
$ws.Cells.Item(53,9).Value = 100101001
$ws.Cells.Item(53,10).Value = "Arándano (blue)"
$ws.Cells.Item(53,11).Value = "Sin especificar"
$ws.Cells.Item(53,12).Value = "Segunda"
$ws.Cells.Item(53,13).Value = 100
$ws.Cells.Item(53,14).Value = 2800
$ws.Cells.Item(53,15).Value = 2800
$ws.Cells.Item(53,16).Value = 2800
$ws.Cells.Item(53,17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(53,18).Value = "Región de Ñuble"
$ws.Cells.Item(53,19).Value = 1400
$ws.Cells.Item(53,20).Value = 2
